$wb = $excel.ActiveWorkbook

# Update the header in the "Expense" sheet: "Latest Payment Date" -> "Payment Date"
$expenseSheet = $wb.Worksheets.Item("Expense")
$expenseSheet.Range("E1").Value = "Payment Date"

# Switch the active sheet to "Expense" and set its selection to E1
$expenseSheet.Activate()
$expenseSheet.Range("E1").Select()

# Keep the "Funding" sheet's prior selection as-is (H7), no changes needed there
